$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version, Date, Publisher, replace duplicate
#     Contact rows with a single Jurisdiction row ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date updated
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes
# "Jurisdiction" / "United States of America"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# it is removed entirely, shifting everything below up by one row.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements": update the Short / Definition text for the root
#     Extension element (row 2) to match the new Title / Description ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Cells.Item(2, 11).Value = "Episode Allowed Amount PCP (USD)"
$elements.Cells.Item(2, 12).Value = "Allowed amount per episide for the primary care physician, in USD"
